$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 - this shifts the existing rows 12-24
# down to 13-25 (carrying their data/formatting with them), matching the
# diff where every row from 12 downward is pushed down by one and a new
# data row appears at row 12.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new weekly price record.
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C12").Value = 'Los Lagos'
$ws.Range("D12").Value = 44799
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 100112013
$ws.Range("G12").Value = 'Alcachofa'
$ws.Range("H12").Value = 'Madrigal'
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 70
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("N12").Value = '$/caja 40 unidades'
$ws.Range("O12").Value = 'Provincia de Limarí'
$ws.Range("P12").Value = 375
$ws.Range("Q12").Value = 40
$ws.Range("R12").Value = 'Hortaliza'
